$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The uploaded-file crash fix: the header cell A1 ("SKU") is replaced with "VPN"
$ws.Range("A1").Value = "VPN"

# Move the active selection to A6 (as reflected in the saved view state)
$ws.Range("A6").Select()
